$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New entry being handed back: ca0d5071-346c-4d2d-8615-3b791614ff99.md
# ---------------------------------------------------------------------------
$fileId = "ca0d5071-346c-4d2d-8615-3b791614ff99"
$mdName = "$fileId.md"
$pathAndName = "e2e\$fileId.md"
$zhXlf = "$fileId.6baa2f45d61d806c03a78ad3f557182a0e89ef7e.zh-cn.xlf"
$deXlf = "$fileId.6baa2f45d61d806c03a78ad3f557182a0e89ef7e.de-de.xlf"

$hoDateZh = "2016-08-31 02:50:37"
$hoDateDe = "2016-08-31 02:50:42"
$overviewDate = "2016-08-31 02:50:42"

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca0d5071346c4d2d86153b791614ff990000000/e2e/$mdName"
$ovUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca0d5071346c4d2d86153b791614ff990000000/e2e/$mdName"

$blueColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 9
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A9").Value = $mdName
$ws1.Range("C9").Value = ".md"
$ws1.Range("D9").Value = ""
$ws1.Range("E9").Value = "Ready for handoff"
$ws1.Range("F9").Value = "Ready for handoff"
$ws1.Range("G9").Value = $overviewDate
$ws1.Range("G9").NumberFormat = $dateFormat

$ws1.Hyperlinks.Add($ws1.Range("B9"), $ovUrl, $null, $null, $pathAndName)
$ws1.Range("B9").Font.Underline = $true
$ws1.Range("B9").Font.Color = $blueColor

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G9"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 9
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Add($ws2.Range("A9"), $zhUrl, $null, $null, $mdName)
$ws2.Range("A9").Font.Underline = $true
$ws2.Range("A9").Font.Color = $blueColor

$ws2.Range("B9").Value = ".md"
$ws2.Range("C9").Value = "Ready for handoff"
$ws2.Range("D9").Value = "e2e"
$ws2.Range("E9").Value = "ht"
$ws2.Range("F9").Value = "False"
$ws2.Range("G9").Value = $zhXlf
$ws2.Range("H9").Value = $hoDateZh
$ws2.Range("H9").NumberFormat = $dateFormat
$ws2.Range("I9").Value = ""
$ws2.Range("J9").Value = ""
$ws2.Range("K9").Value = "0001-01-01 00:00:00"
$ws2.Range("K9").NumberFormat = $dateFormat
$ws2.Range("L9").Value = ""
$ws2.Range("M9").Value = "True"
$ws2.Range("N9").Value = ""
$ws2.Range("O9").Value = "False"
$ws2.Range("P9").Value = ""

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P9"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 9
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Add($ws3.Range("A9"), $zhUrl, $null, $null, $mdName)
$ws3.Range("A9").Font.Underline = $true
$ws3.Range("A9").Font.Color = $blueColor

$ws3.Range("B9").Value = ".md"
$ws3.Range("C9").Value = "Ready for handoff"
$ws3.Range("D9").Value = "e2e"
$ws3.Range("E9").Value = "ht"
$ws3.Range("F9").Value = "False"
$ws3.Range("G9").Value = $deXlf
$ws3.Range("H9").Value = $hoDateDe
$ws3.Range("H9").NumberFormat = $dateFormat
$ws3.Range("I9").Value = ""
$ws3.Range("J9").Value = ""
$ws3.Range("K9").Value = "0001-01-01 00:00:00"
$ws3.Range("K9").NumberFormat = $dateFormat
$ws3.Range("L9").Value = ""
$ws3.Range("M9").Value = "True"
$ws3.Range("N9").Value = ""
$ws3.Range("O9").Value = "False"
$ws3.Range("P9").Value = ""

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P9"))
